$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the steel description text in B2: remove the "/RME" segment
# from the second bullet line ("25% S/LFM+CDN/RME/H:1" -> "25% S/LFM+CDN/H:1")
$newText = "20% S+SL/LFM+CDN/H:1`n25% S/LFM+CDN/H:1`n25% CR+PC/LFM+CDN/H:1`n10% CR/LWAL+CDN/H:2`n8% CR/LFM+CDN/H:1`n5% W/LWAL+CDN/H:1`n7% MUR/LWAL+CDN/H:1"
$ws.Range("B2").Value = $newText

# Wrap the text in B2 and resize the row to fit the multi-line content
$ws.Range("B2").WrapText = $true
$ws.Rows.Item(2).RowHeight = 365

# Update the active selection in the sheet view
$ws.Range("B9").Select()
